# Append the next email address to the "Emails" sheet.
# The existing data lives in column A, rows 1-10; this adds row 11,
# growing the used range (and its dimension) to A1:A11.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "abibangbrandon855@gmail.com"
